$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds numeric-looking text (e.g. "1.00", "0.0980",
# "57.572.14" using dots as thousand separators). Setting NumberFormat to
# "@" (text) before writing the value keeps it verbatim instead of letting
# Excel auto-convert/round it as a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.572.14"
$ws.Range("E2").Value = "  -5.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.895.89"
$ws.Range("E3").Value = "  -3.92%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.45"
$ws.Range("E5").Value = "  -3.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.16"
$ws.Range("E6").Value = "  -4.40%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.890.00"
$ws.Range("E8").Value = "  -4.08%  "
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  -7.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.73"
$ws.Range("E11").Value = "  -9.00%  "
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("E13").Value = "  -4.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.31"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.370.09"
$ws.Range("E16").Value = "  -4.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.894.52"
$ws.Range("E17").Value = "  -3.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.57"
$ws.Range("E18").Value = "  +5.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.496.62"
$ws.Range("E19").Value = "  -6.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "407.27"
$ws.Range("E20").Value = "  -7.14%  "
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.83"
$ws.Range("E23").Value = "  -4.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.86"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "77.05"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.46"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("E29").Value = "  +3.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.21"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.02"
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.65"
$ws.Range("E32").Value = "  -3.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0980"
$ws.Range("E33").Value = "  +4.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.911"
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.01"
$ws.Range("E36").Value = "  -11.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.17"
$ws.Range("E37").Value = "  -3.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.46"
$ws.Range("E38").Value = "  +9.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0619"
$ws.Range("E39").Value = "  -8.34%  "
$ws.Range("E40").Value = "  -5.63%  "
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.614.36"
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "360.35"
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.41"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "120.10"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  -3.05%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.85"
$ws.Range("E50").Value = "  -2.98%  "
$ws.Range("E51").Value = "  -3.81%  "
